$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row of raw/clean SSA data for Jul 3rd at row 34
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "2020-07-03"
$ws.Range("A34").Style = "Normal"
$ws.Range("B34").Value = 245251
$ws.Range("C34").Value = 301986
$ws.Range("D34").Value = 77750
$ws.Range("E34").Value = 29843
$ws.Range("F34").Value = 30.35
